$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the two new audit-log rows.
#    Row 12: "Failed to Create User (Invalid Password)" - right after
#    "Failed to Create User (Unauthorized)".
# ---------------------------------------------------------------------
$ws.Rows(12).Insert()
$ws.Range("B12").Value = "Failed to Create User (Invalid Password)"
$ws.Range("C12").Value = "Yes"

# Row 15 (after the first insert shifted things down): "Failed to Update
# User (Invalid Password)" - right after "Failed to Update User
# (Unauthorized)" (now at row 14).
$ws.Rows(15).Insert()
$ws.Range("B15").Value = "Failed to Update User (Invalid Password)"
$ws.Range("C15").Value = "Yes"

# ---------------------------------------------------------------------
# 2. Prepare a clean "style 4" (new alternating-row fill: theme 0 / white,
#    no tint) in a scratch cell far away from the used range. We seed it
#    from an already-explicitly-styled cell (B4, style 2) so the engine's
#    fill writer takes the in-place-replace path instead of allocating an
#    orphan intermediate fill (which happens when writing to a default
#    -styled cell). This runs after the row inserts so the scratch cell
#    location is not disturbed by further row shifting.
# ---------------------------------------------------------------------
$ws.Range("B4").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("Z100").Interior.ThemeColor = 2
$ws.Range("Z100").Interior.TintAndShade = 0

# ---------------------------------------------------------------------
# 3. Re-apply the alternating row banding from row 12 down to row 20
#    (style 3 = existing grey fill [formerly only on "View Audit Logs",
#    now at row 18], style 4 = the new white fill prepared above).
# ---------------------------------------------------------------------
$ws.Range("B18").Copy()
$ws.Range("B12:C12").PasteSpecial(-4122)
$ws.Range("B14:C14").PasteSpecial(-4122)
$ws.Range("B16:C16").PasteSpecial(-4122)
$ws.Range("B18:C18").PasteSpecial(-4122)
$ws.Range("B20:C20").PasteSpecial(-4122)

$ws.Range("Z100").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Range("B15:C15").PasteSpecial(-4122)
$ws.Range("B17:C17").PasteSpecial(-4122)
$ws.Range("B19:C19").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4. Clean up the scratch cell so it doesn't leak into the used range.
# ---------------------------------------------------------------------
$ws.Range("Z100").Clear()

# ---------------------------------------------------------------------
# 5. Restore selection to match the post-edit workbook (B11).
# ---------------------------------------------------------------------
$ws.Range("B11").Select()
